# Update "想去人数" (F column) values across the "展览", "演出" and "全部类型"
# worksheets to reflect the latest scrape, as published to gh-pages.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 13779
    3  = 89
    4  = 127
    5  = 545
    7  = 1205
    8  = 1028
    9  = 13890
    10 = 14743
    12 = 4
    21 = 61
    23 = 1148
    26 = 5719
    27 = 943
    28 = 1056
    29 = 5406
    30 = 46
    32 = 255
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 6

# --- Sheet: 全部类型 (All types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 13779
    3  = 89
    4  = 6
    5  = 127
    6  = 545
    8  = 1205
    9  = 1028
    10 = 13890
    11 = 14743
    13 = 4
    22 = 61
    24 = 1148
    27 = 5719
    28 = 943
    29 = 1056
    30 = 5406
    31 = 46
    33 = 255
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
